# Fix bugs, added review, results
# The "correct" column (H) for true_false questions incorrectly stored the
# literal string "true" instead of referencing the actual correct option
# (option_1 / "Option A"). Update every such row to say "option_1".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 4, 6, 8, 10, 12, 14, 16)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 8).Value = "option_1"
}

# Restore the active selection left behind by the review pass.
$ws.Range("K10").Select() | Out-Null
